$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Re-assert the original sheet's column widths (13 characters) so the
# "true"/"false"-style customWidth flag round-trips correctly.
$ws1.Range("A1:K1").EntireColumn.ColumnWidth = 12.166666666666666

# Add a new worksheet after the existing one and name it "Sheet1"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

# Copy the daily-data range (header + 31 days of data) from the original
# sheet into the new sheet, starting at A1. Two PasteSpecial passes so we
# carry over both the cell formatting (borders/alignment) and the values.
$src = $ws1.Range("A9:K40")
$src.Copy()
$dst = $ws2.Range("A1")
$dst.PasteSpecial(-4122)
$dst.PasteSpecial(-4163)

# Approximate Excel's auto-fit row height (wrapped text in the narrower
# default columns of the new sheet makes each data row take two lines).
$ws2.Range("A2:K32").RowHeight = 28.8

# Update the view/selection on the original sheet.
$null = $ws1.Range("A9:K40").Select()

# Make the new sheet the active tab with the pasted range selected.
$null = $ws2.Range("A1:K32").Select()
$null = $ws2.Activate()
